# Chelsea_stats.xlsx update
# 1) Update a handful of "age" values (format YY-DDD, day-of-year) in column D
#    of every stats sheet (Standard Stats .. Miscellaneous Stats) - the data
#    was refreshed a few days later, so each age/day figure advances by 3 days
#    (with year rollover when the day count exceeds 364).
# 2) Fix the Date/Time/Day of the Liverpool fixture (row 52) on the Matches
#    sheet (moved from Sat 2025-05-03 15:00 to Sun 2025-05-04 16:30).
# 3) Shift every sheet name along by one position (Matches is dropped from
#    the tab list, every other tab takes over the name of the sheet that used
#    to precede it) and the last sheet becomes "Sheet_9".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data: the 29 "age" values shared by every stats sheet (rows 4-32),
# plus the extra 6 rows (33-38) that only exist on "Standard Stats" and
# "Playing Time".
# ---------------------------------------------------------------------------
$oldBase = @("23-142","22-322","22-026","26-245","24-066","27-126","23-277","23-014","21-309","25-015","24-364","24-097","27-181","27-130","25-106","21-077","22-342","27-013","25-262","25-134","23-363","18-323","26-199","21-238","24-078","19-048","19-079","19-065","17-155")
$newBase = @("23-145","22-325","22-029","26-248","24-069","27-129","23-280","23-017","21-312","25-018","25-002","24-100","27-184","27-133","25-109","21-080","22-345","27-016","25-265","25-137","24-001","18-326","26-202","21-241","24-081","19-051","19-082","19-068","17-158")

$oldExtra = @("19-329","17-317","32-304","22-073","19-362","18-292")
$newExtra = @("19-332","17-320","32-307","22-076","20-000","18-295")

# Sheet names as they currently exist (before the rename shuffle below).
$statSheets29 = @("Shooting Stats","Passing Stats","Pass Types","Goal & Shot Creation","Defensive Actions","Possession","Miscellaneous Stats")
$statSheets35 = @("Standard Stats","Playing Time")

function Update-AgeColumn($sheetName, $startRow, $oldVals, $newVals) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $newVals.Length; $i++) {
        $row = $startRow + $i
        # "YY-DDD" strings are never mistaken for dates/numbers by Excel, so
        # a plain value assignment keeps them stored as text, same as before.
        $ws.Cells.Item($row, 4).Value = $newVals[$i]
    }
}

foreach ($name in $statSheets29) {
    Update-AgeColumn $name 4 $oldBase $newBase
}
foreach ($name in $statSheets35) {
    Update-AgeColumn $name 4 $oldBase $newBase
    Update-AgeColumn $name 33 $oldExtra $newExtra
}

# ---------------------------------------------------------------------------
# Matches sheet: correct the Liverpool match on row 52.
# ---------------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches")

# "2025-05-04" reads as an ISO date, so Excel's COM layer would silently
# convert it to a date serial number; force the cell to Text first, assign
# the value, then restore the original ("Normal"/General) cell style so the
# only semantic change is the text itself.
$matches.Cells.Item(52, 2).NumberFormat = "@"
$matches.Cells.Item(52, 2).Value = "2025-05-04"
$matches.Cells.Item(52, 2).Style = "Normal"

# "16:30" is not auto-converted like the date above, so a plain assignment
# is enough to keep it stored as text.
$matches.Cells.Item(52, 3).Value = "16:30"

$matches.Cells.Item(52, 6).Value = "Sun"

# ---------------------------------------------------------------------------
# Shuffle sheet names: drop "Matches" from the tab list and shift every
# remaining tab name back one slot, appending "Sheet_9" at the end.
# Rename from the last sheet to the first to avoid name collisions.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Miscellaneous Stats").Name = "Sheet_9"
$wb.Worksheets.Item("Playing Time").Name = "Miscellaneous Stats"
$wb.Worksheets.Item("Possession").Name = "Playing Time"
$wb.Worksheets.Item("Defensive Actions").Name = "Possession"
$wb.Worksheets.Item("Goal & Shot Creation").Name = "Defensive Actions"
$wb.Worksheets.Item("Pass Types").Name = "Goal & Shot Creation"
$wb.Worksheets.Item("Passing Stats").Name = "Pass Types"
$wb.Worksheets.Item("Shooting Stats").Name = "Passing Stats"
$wb.Worksheets.Item("Standard Stats").Name = "Shooting Stats"
$wb.Worksheets.Item("Matches").Name = "Standard Stats"
